# Fruta / hortaliza, semanal
# Re-shuffle the weekly Fecha/Volumen/Precio* figures (cols D, J, K, L, M, P)
# across rows 2-31 of the Ciboulette sheet -- each row keeps its own
# market/category/quality metadata but picks up another row's
# date + volume + price-range values (weekly re-roll of the sample).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2=44978, J2=1000, K2=1800, L2=2000, M2=1900, P2=633
$ws.Range("D2").Value = 44978
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1800
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1900
$ws.Range("P2").Value = 633

# Row 3: D3=44965, J3=1120, K3=2000, L3=2500, M3=2250, P3=750
$ws.Range("D3").Value = 44965
$ws.Range("J3").Value = 1120
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2250
$ws.Range("P3").Value = 750

# Row 4: D4=45006, J4=1100
$ws.Range("D4").Value = 45006
$ws.Range("J4").Value = 1100

# Row 5: D5=44848, K5=1500, L5=2000, M5=1750, P5=583
$ws.Range("D5").Value = 44848
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1750
$ws.Range("P5").Value = 583

# Row 6: D6=45020, J6=1200, K6=2000, L6=2500, M6=2250, P6=750
$ws.Range("D6").Value = 45020
$ws.Range("J6").Value = 1200
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = 2250
$ws.Range("P6").Value = 750

# Row 7: D7=44985
$ws.Range("D7").Value = 44985

# Row 8: D8=44881, J8=500, K8=1900, M8=1950, P8=650
$ws.Range("D8").Value = 44881
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 1900
$ws.Range("M8").Value = 1950
$ws.Range("P8").Value = 650

# Row 9: D9=45070, J9=800
$ws.Range("D9").Value = 45070
$ws.Range("J9").Value = 800

# Row 10: D10=44999, J10=1100, K10=2000, L10=2500, M10=2250, P10=750
$ws.Range("D10").Value = 44999
$ws.Range("J10").Value = 1100
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 2500
$ws.Range("M10").Value = 2250
$ws.Range("P10").Value = 750

# Row 11: D11=44971, J11=1000
$ws.Range("D11").Value = 44971
$ws.Range("J11").Value = 1000

# Row 12: D12=44992, J12=1040
$ws.Range("D12").Value = 44992
$ws.Range("J12").Value = 1040

# Row 13: D13=45062
$ws.Range("D13").Value = 45062

# Row 14: D14=44827, J14=1200, K14=2000, L14=2500, M14=2250, P14=750
$ws.Range("D14").Value = 44827
$ws.Range("J14").Value = 1200
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2250
$ws.Range("P14").Value = 750

# Row 15: D15=44911, J15=700, K15=1800, L15=2000, M15=1900, P15=633
$ws.Range("D15").Value = 44911
$ws.Range("J15").Value = 700
$ws.Range("K15").Value = 1800
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 1900
$ws.Range("P15").Value = 633

# Row 16: D16=44970, J16=800
$ws.Range("D16").Value = 44970
$ws.Range("J16").Value = 800

# Row 17: D17=45035, J17=1100
$ws.Range("D17").Value = 45035
$ws.Range("J17").Value = 1100

# Row 18: D18=45034, J18=1100
$ws.Range("D18").Value = 45034
$ws.Range("J18").Value = 1100

# Row 19: D19=44964, J19=1000
$ws.Range("D19").Value = 44964
$ws.Range("J19").Value = 1000

# Row 20: D20=45041, J20=1160, K20=2000, L20=2500, M20=2250, P20=750
$ws.Range("D20").Value = 45041
$ws.Range("J20").Value = 1160
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = 2250
$ws.Range("P20").Value = 750

# Row 21: D21=45028, J21=1000
$ws.Range("D21").Value = 45028
$ws.Range("J21").Value = 1000

# Row 22: D22=45013, J22=1100, K22=2000, L22=2500, M22=2250, P22=750
$ws.Range("D22").Value = 45013
$ws.Range("J22").Value = 1100
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = 2250
$ws.Range("P22").Value = 750

# Row 23: D23=44685, J23=400, K23=1500, L23=2000, M23=1750, P23=583
$ws.Range("D23").Value = 44685
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 2000
$ws.Range("M23").Value = 1750
$ws.Range("P23").Value = 583

# Row 24: D24=45084, J24=900
$ws.Range("D24").Value = 45084
$ws.Range("J24").Value = 900

# Row 25: D25=44883, J25=500, K25=1800, L25=2000, M25=1900, P25=633
$ws.Range("D25").Value = 44883
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 1800
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 1900
$ws.Range("P25").Value = 633

# Row 26: D26=44951, J26=800
$ws.Range("D26").Value = 44951
$ws.Range("J26").Value = 800

# Row 27: D27=45091, J27=800
$ws.Range("D27").Value = 45091
$ws.Range("J27").Value = 800

# Row 28: D28=45007, J28=1160
$ws.Range("D28").Value = 45007
$ws.Range("J28").Value = 1160

# Row 29: D29=44910, J29=1000, K29=1800, L29=2000, M29=1900, P29=633
$ws.Range("D29").Value = 44910
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 1800
$ws.Range("L29").Value = 2000
$ws.Range("M29").Value = 1900
$ws.Range("P29").Value = 633

# Row 30: D30=45077, J30=760
$ws.Range("D30").Value = 45077
$ws.Range("J30").Value = 760

# Row 31: D31=44953
$ws.Range("D31").Value = 44953
